# Update the "想去人数" (want-to-go count) figures that changed between
# this gh-pages data refresh and the previous one.
#
# Both the "展览" sheet and the "全部类型" sheet contain the same rows of
# exhibition data, so the same three cells need to be bumped in each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F13" = 1232
    "F15" = 381
    "F22" = 2707
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
